# Update the attendance sheet for tut06 / 2001EE56.xlsx
# - Change date strings in column A from DD/MM/YYYY to DD-MM-YYYY format (rows 3-21)
# - Update some D/E/G/H numeric values on specific rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (slashes replaced with dashes)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Rows whose day-of-month is <= 12, so Excel's automatic type detection
# could otherwise misinterpret "DD-MM-YYYY" as a date value. For those,
# force the cell to Text format first so the literal string is kept.
$ambiguousRows = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    if ($ambiguousRows -contains $row) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $dates[$row]
}

# Numeric cell updates: row -> column -> new value
$updates = @{
    3  = @{ D = 1; G = 1 }
    5  = @{ D = 1; E = 1; H = 0 }
    11 = @{ D = 1; E = 1; H = 0 }
    12 = @{ D = 1; E = 1; H = 0 }
    13 = @{ D = 1; E = 1; H = 0 }
    15 = @{ D = 1; E = 1; H = 0 }
    16 = @{ D = 1; E = 1; H = 0 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
